$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "Modelo"

$ws.Range("F2").Value = "Pipeline(steps=[('model', LinearRegression())])"
